$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- New row 13: Agarwal 2020 - substantia nigra ---
$ws.Range("A13").Value = "Agarwal2020_SNIG"
$ws.Range("B13").Value = "Agarwal"
$ws.Range("C13").Value = 2020
$ws.Range("D13").Value = "https://www.ncbi.nlm.nih.gov/geo/query/acc.cgi?acc=GSE140231"
$ws.Range("E13").Value = "human"
$ws.Range("F13").Value = "substantia_nigra"
$ws.Range("G13").Value = "control;post_mortem"
$ws.Range("H13").Value = "single_nuclei"
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 6105
$ws.Range("K13").Value = 7
$ws.Range("L13").Value = 10
$ws.Range("M13").Value = "https://www.ncbi.nlm.nih.gov/geo/query/acc.cgi?acc=GSE140231"
$ws.Range("N13").Value = "GSE140231_RAW.tar"
$ws.Range("O13").Value = "SI_matrix"

# --- New row 14: Agarwal 2020 - cortex ---
$ws.Range("A14").Value = "Agarwal2020_CRTX"
$ws.Range("B14").Value = "Agarwal"
$ws.Range("C14").Value = 2020
$ws.Range("D14").Value = "https://www.ncbi.nlm.nih.gov/geo/query/acc.cgi?acc=GSE140231"
$ws.Range("E14").Value = "human"
$ws.Range("F14").Value = "cortex"
$ws.Range("G14").Value = "control;post_mortem"
$ws.Range("H14").Value = "single_nuclei"
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 12015
$ws.Range("K14").Value = 6
$ws.Range("L14").Value = 23
$ws.Range("M14").Value = "https://www.ncbi.nlm.nih.gov/geo/query/acc.cgi?acc=GSE140231"
$ws.Range("N14").Value = "GSE140231_RAW.tar"
$ws.Range("O14").Value = "SI_matrix"

# Match the formatting used on the row above (N column carries a "no fill"
# style flag in the source row; re-apply the same formatting to the new cells)
$ws.Range("N12").Copy()
$ws.Range("N13:N14").PasteSpecial(-4122)
$ws.Range("N13").Value = "GSE140231_RAW.tar"
$ws.Range("N14").Value = "GSE140231_RAW.tar"

# View changes: zoom back to 100% and move the active selection past the new rows
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("A15").Select()
